$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 8230830
$ws.Range("J17").Value = 8230830
$ws.Range("L17").Value = 24692490
$ws.Range("N17").Value = -24692826

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 78999.5
$ws.Range("J120").Value = 78999.5
$ws.Range("L120").Value = 78999.5
$ws.Range("N120").Value = -88675.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2331.4773
$ws.Range("I132").Value = 2142
$ws.Range("K132").Value = 6426
$ws.Range("M132").Value = -3896

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 27856.908
$ws.Range("I137").Value = 57467
$ws.Range("K137").Value = 172401
$ws.Range("M137").Value = -169851

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3196.5
$ws.Range("I138").Value = 3196.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9589.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -4449.5
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4961.9263
$ws.Range("J32").Value = 33604
$ws.Range("L32").Value = 33604
$ws.Range("N32").Value = -34178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 95031
$ws.Range("J42").Value = 95031
$ws.Range("L42").Value = 95031
$ws.Range("N42").Value = -96003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4635.25
$ws.Range("I45").Value = 4635.25
$ws.Range("K45").Value = 4635.25
$ws.Range("M45").Value = -4258.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3249.125
$ws.Range("I63").Value = 3599
$ws.Range("J63").Value = 2666
$ws.Range("K63").Value = 3599
$ws.Range("L63").Value = 2666
$ws.Range("M63").Value = -2913
$ws.Range("N63").Value = -4038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3249.125
$ws.Range("I66").Value = 3599
$ws.Range("J66").Value = 2666
$ws.Range("K66").Value = 17995
$ws.Range("L66").Value = 13330
$ws.Range("M66").Value = -14563
$ws.Range("N66").Value = -20194

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2577.6956
$ws.Range("I132").Value = 2514.5
$ws.Range("K132").Value = 7543.5
$ws.Range("M132").Value = -5013.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 733.3333
$ws.Range("I22").Value = 731.8182
$ws.Range("K22").Value = 731.8182
$ws.Range("M22").Value = -558.8182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 863.9231
$ws.Range("I80").Value = 1080.9375
$ws.Range("J80").Value = 516.7
$ws.Range("K80").Value = 1080.9375
$ws.Range("L80").Value = 516.7
$ws.Range("M80").Value = -82.9375
$ws.Range("N80").Value = -2512.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 863.9231
$ws.Range("I83").Value = 1080.9375
$ws.Range("J83").Value = 516.7
$ws.Range("K83").Value = 5404.6875
$ws.Range("L83").Value = 2583.5
$ws.Range("M83").Value = -412.6875
$ws.Range("N83").Value = -12567.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3466.9834
$ws.Range("I134").Value = 3594.698
$ws.Range("K134").Value = 10784.094
$ws.Range("M134").Value = -8249.093999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 135660.05
$ws.Range("I31").Value = 198104.8
$ws.Range("J31").Value = 2964.9583
$ws.Range("K31").Value = 198104.8
$ws.Range("L31").Value = 2964.9583
$ws.Range("M31").Value = -197809.8
$ws.Range("N31").Value = -3554.9583

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 135660.05
$ws.Range("I34").Value = 198104.8
$ws.Range("J34").Value = 2964.9583
$ws.Range("K34").Value = 198104.8
$ws.Range("L34").Value = 2964.9583
$ws.Range("M34").Value = -197902.8
$ws.Range("N34").Value = -3368.9583

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2718.5667
$ws.Range("J58").Value = 3228.8
$ws.Range("L58").Value = 3228.8
$ws.Range("N58").Value = -3634.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2132.2
$ws.Range("I134").Value = 2162.0256
$ws.Range("J134").Value = 969
$ws.Range("K134").Value = 6486.0768
$ws.Range("L134").Value = 2907
$ws.Range("M134").Value = -3951.0768
$ws.Range("N134").Value = -7977

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2718.5667
$ws.Range("J136").Value = 3228.8
$ws.Range("L136").Value = 9686.400000000001
$ws.Range("N136").Value = -14786.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 562.8570999999999
$ws.Range("I5").Value = 608.2
$ws.Range("J5").Value = 449.5
$ws.Range("K5").Value = 1824.6
$ws.Range("L5").Value = 1348.5
$ws.Range("M5").Value = -1712.6
$ws.Range("N5").Value = -1572.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 3395
$ws.Range("J32").Value = 3994
$ws.Range("L32").Value = 11982
$ws.Range("N32").Value = -12548

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1111.5
$ws.Range("J46").Value = 1473.75
$ws.Range("L46").Value = 4421.25
$ws.Range("N46").Value = -4603.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2607537.5
$ws.Range("I68").Value = 8335051
$ws.Range("J68").Value = 4122.1816
$ws.Range("K68").Value = 25005153
$ws.Range("L68").Value = 12366.5448
$ws.Range("M68").Value = -25004342
$ws.Range("N68").Value = -13988.5448

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2607537.5
$ws.Range("I71").Value = 8335051
$ws.Range("J71").Value = 4122.1816
$ws.Range("K71").Value = 75015459
$ws.Range("L71").Value = 37099.6344
$ws.Range("M71").Value = -75011403
$ws.Range("N71").Value = -45211.6344

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 2105.125
$ws.Range("J112").Value = 2338
$ws.Range("L112").Value = 7014
$ws.Range("N112").Value = -9230

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 125128690
$ws.Range("I121").Value = 250002500
$ws.Range("J121").Value = 254873.75
$ws.Range("K121").Value = 750007500
$ws.Range("L121").Value = 764621.25
$ws.Range("M121").Value = -750006190
$ws.Range("N121").Value = -767241.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 562.8570999999999
$ws.Range("I135").Value = 608.2
$ws.Range("J135").Value = 449.5
$ws.Range("K135").Value = 5473.8
$ws.Range("L135").Value = 4045.5
$ws.Range("M135").Value = -2938.8
$ws.Range("N135").Value = -9115.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2363.55
$ws.Range("I139").Value = 1318.0667
$ws.Range("J139").Value = 5500
$ws.Range("K139").Value = 3954.2001
$ws.Range("L139").Value = 16500
$ws.Range("M139").Value = 1185.7999
$ws.Range("N139").Value = -26780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2866.9
$ws.Range("I70").Value = 2464.6924
$ws.Range("J70").Value = 3613.8572
$ws.Range("K70").Value = 2464.6924
$ws.Range("L70").Value = 3613.8572
$ws.Range("M70").Value = -2194.6924
$ws.Range("N70").Value = -4153.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 2866.9
$ws.Range("I73").Value = 2464.6924
$ws.Range("J73").Value = 3613.8572
$ws.Range("K73").Value = 2464.6924
$ws.Range("L73").Value = 3613.8572
$ws.Range("M73").Value = -1528.6924
$ws.Range("N73").Value = -5485.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2904.8235
$ws.Range("I97").Value = 2892.1333
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 2892.1333
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -2396.1333
$ws.Range("N97").Value = -3992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 30560.441
$ws.Range("I107").Value = 44621.78
$ws.Range("K107").Value = 44621.78
$ws.Range("M107").Value = -42701.78

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 14461.5
$ws.Range("I126").Value = 14695.579
$ws.Range("J126").Value = 10014
$ws.Range("K126").Value = 44086.737
$ws.Range("L126").Value = 30042
$ws.Range("M126").Value = -41616.737
$ws.Range("N126").Value = -34982

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 23244.51
$ws.Range("I132").Value = 24186.852
$ws.Range("J132").Value = 1099.5
$ws.Range("K132").Value = 72560.556
$ws.Range("L132").Value = 3298.5
$ws.Range("M132").Value = -70030.556
$ws.Range("N132").Value = -8358.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6148.4287
$ws.Range("I7").Value = 6000.5
$ws.Range("K7").Value = 6000.5
$ws.Range("M7").Value = -5888.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 35716310
$ws.Range("I16").Value = 41668484
$ws.Range("J16").Value = 3300.5
$ws.Range("K16").Value = 41668484
$ws.Range("L16").Value = 3300.5
$ws.Range("M16").Value = -41668314
$ws.Range("N16").Value = -3640.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6148.4287
$ws.Range("I126").Value = 6000.5
$ws.Range("K126").Value = 18001.5
$ws.Range("M126").Value = -15531.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3782.2693
$ws.Range("I136").Value = 3223.4348
$ws.Range("K136").Value = 9670.304400000001
$ws.Range("M136").Value = -7120.304400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1333.76
$ws.Range("I122").Value = 1276.3414
$ws.Range("K122").Value = 3829.0242
$ws.Range("M122").Value = -1379.0242

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 65471.477
$ws.Range("I126").Value = 2378.9412
$ws.Range("J126").Value = 244233.67
$ws.Range("K126").Value = 7136.823600000001
$ws.Range("L126").Value = 732701.01
$ws.Range("M126").Value = -4666.823600000001
$ws.Range("N126").Value = -737641.01

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1910.7391
$ws.Range("I132").Value = 927.55
$ws.Range("J132").Value = 2667.0386
$ws.Range("K132").Value = 2782.65
$ws.Range("L132").Value = 8001.1158
$ws.Range("M132").Value = -252.6499999999996
$ws.Range("N132").Value = -13061.1158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 436385.97
$ws.Range("I136").Value = 589347.0600000001
$ws.Range("J136").Value = 2996.1667
$ws.Range("K136").Value = 1768041.18
$ws.Range("L136").Value = 8988.500100000001
$ws.Range("M136").Value = -1765491.18
$ws.Range("N136").Value = -14088.5001
